$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update summary figures at top of the statement ---
# Valor Mora (total past-due amount)
$ws.Range("E11").Value = 255692
# Cant. Trabajadores (worker count)
$ws.Range("C13").Value = 4
# Cant. Periodos (period count)
$ws.Range("F13").Value = 2

# --- Update the existing detail row (was LAURIANO CURE SUAREZ / 1046404907) ---
# to reflect the refreshed record for YORMAN SANTIAGO AYALA / 20246181
$ws.Range("C17").Value = "20246181"
$ws.Range("D17").Value = "YORMAN SANTIAGO AYALA"
$ws.Range("F17").Value = 33120
$ws.Range("G17").Value = 828000

# --- Insert a new detail row right after row 17, pushing the former last row   ---
# --- (the bottom-bordered YORMAN row, now stale) down to become the new row 19 ---
$ws.Rows("18:18").Insert()

# Copy the border/fill/font formatting from row 17 onto the freshly inserted row 18
$ws.Range("B17:J17").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)

# Populate the new row 18 with the first added worker (part 1 of new statements)
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1047427541"
$ws.Range("D18").Value = "FABIO EDUARDO GRONDONA PATERNINA"
$ws.Range("E18").Value = "2508"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

# Row 19 is the row that was pushed down; it keeps the bottom-border style and now
# holds the second added worker
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "72166473"
$ws.Range("D19").Value = "JORGE LUIS GRONDONA VILLEGAS"
$ws.Range("E19").Value = "2508"
$ws.Range("F19").Value = 156800
$ws.Range("G19").Value = 3920000
